$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: A2 = "PARA01", B2 = 8801260607365 (numeric barcode)
$ws.Range("A2").Value = "PARA01"
$ws.Range("B2").Value = 8801260607365

# Apply a number format to column B (integer with trailing space)
$ws.Range("B1:B2").NumberFormat = "0_ "

# Autofit column B to match bestFit width behavior
$ws.Columns("B").AutoFit()

# Update selection to mimic final cursor position
$ws.Range("F17").Select()
